$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2391.8333
$ws.Range("I4").Value = 785
$ws.Range("K4").Value = 785
$ws.Range("M4").Value = -671
$ws.Range("H9").Value = 737.3333
$ws.Range("I9").Value = 356.125
$ws.Range("K9").Value = 356.125
$ws.Range("M9").Value = -187.125
$ws.Range("H40").Value = 3622.1853
$ws.Range("I40").Value = 3266.3333
$ws.Range("J40").Value = 3800.111
$ws.Range("K40").Value = 3266.3333
$ws.Range("L40").Value = 3800.111
$ws.Range("M40").Value = -3091.3333
$ws.Range("N40").Value = -4150.111
$ws.Range("H44").Value = 8749.799999999999
$ws.Range("J44").Value = 34999
$ws.Range("L44").Value = 34999
$ws.Range("N44").Value = -35923
$ws.Range("H70").Value = 1048.2927
$ws.Range("J70").Value = 946.80554
$ws.Range("L70").Value = 2840.41662
$ws.Range("N70").Value = -3380.41662
$ws.Range("H73").Value = 1048.2927
$ws.Range("J73").Value = 946.80554
$ws.Range("L73").Value = 2840.41662
$ws.Range("N73").Value = -4712.41662
$ws.Range("H113").Value = 7149.25
$ws.Range("I113").Value = 3877.6
$ws.Range("K113").Value = 3877.6
$ws.Range("M113").Value = -623.5999999999999
$ws.Range("H135").Value = 1716.2106
$ws.Range("I135").Value = 493.92307
$ws.Range("K135").Value = 4445.30763
$ws.Range("M135").Value = -1910.30763
$ws.Range("H137").Value = 2804.805
$ws.Range("J137").Value = 3759.7058
$ws.Range("L137").Value = 11279.1174
$ws.Range("N137").Value = -16379.1174
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 26666.777
$ws.Range("I50").Value = 10000.333
$ws.Range("J50").Value = 59999.668
$ws.Range("K50").Value = 10000.333
$ws.Range("L50").Value = 59999.668
$ws.Range("M50").Value = -9375.333000000001
$ws.Range("N50").Value = -61249.668
$ws.Range("H60").Value = 20555.334
$ws.Range("I60").Value = 16499.8
$ws.Range("J60").Value = 25624.75
$ws.Range("K60").Value = 16499.8
$ws.Range("L60").Value = 25624.75
$ws.Range("M60").Value = -15988.8
$ws.Range("N60").Value = -26646.75
$ws.Range("H107").Value = 1092.1936
$ws.Range("I107").Value = 891.53845
$ws.Range("J107").Value = 2135.6
$ws.Range("K107").Value = 891.53845
$ws.Range("L107").Value = 2135.6
$ws.Range("M107").Value = 1028.46155
$ws.Range("N107").Value = -5975.6
$ws.Range("H132").Value = 6518.609
$ws.Range("I132").Value = 3090.0625
$ws.Range("J132").Value = 14355.286
$ws.Range("K132").Value = 9270.1875
$ws.Range("L132").Value = 43065.858
$ws.Range("M132").Value = -6740.1875
$ws.Range("N132").Value = -48125.858
$ws.Range("H134").Value = 6206.0347
$ws.Range("I134").Value = 4452.3184
$ws.Range("J134").Value = 11717.714
$ws.Range("K134").Value = 13356.9552
$ws.Range("L134").Value = 35153.142
$ws.Range("M134").Value = -10821.9552
$ws.Range("N134").Value = -40223.142
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3340
$ws.Range("I3").Value = 3340
$ws.Range("K3").Value = 10020
$ws.Range("M3").Value = -9908
$ws.Range("H8").Value = 417.77777
$ws.Range("I8").Value = 417.77777
$ws.Range("K8").Value = 1253.33331
$ws.Range("M8").Value = -1114.33331
$ws.Range("H98").Value = 917.55554
$ws.Range("I98").Value = 889.75
$ws.Range("J98").Value = 939.8
$ws.Range("K98").Value = 2669.25
$ws.Range("L98").Value = 2819.4
$ws.Range("M98").Value = -1171.25
$ws.Range("N98").Value = -5815.4
$ws.Range("H124").Value = 4027.5
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H133").Value = 13516.583
$ws.Range("I133").Value = 5383.1665
$ws.Range("K133").Value = 16149.4995
$ws.Range("M133").Value = -11089.4995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2859.88
$ws.Range("I102").Value = 2068.1428
$ws.Range("K102").Value = 2068.1428
$ws.Range("M102").Value = -446.1428000000001
$ws.Range("H132").Value = 4130
$ws.Range("I132").Value = 2976.375
$ws.Range("K132").Value = 8929.125
$ws.Range("M132").Value = -6399.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5556.6
$ws.Range("J46").Value = 5556.6
$ws.Range("L46").Value = 5556.6
$ws.Range("N46").Value = -5932.6
$ws.Range("H68").Value = 7583.25
$ws.Range("I68").Value = 7599.6
$ws.Range("K68").Value = 7599.6
$ws.Range("M68").Value = -6850.6
$ws.Range("H71").Value = 7583.25
$ws.Range("I71").Value = 7599.6
$ws.Range("K71").Value = 37998
$ws.Range("M71").Value = -34254
$ws.Range("H100").Value = 2840.2
$ws.Range("I100").Value = 2840.2
$ws.Range("K100").Value = 2840.2
$ws.Range("M100").Value = -2299.2
$ws.Range("H136").Value = 30308550
$ws.Range("I136").Value = 4897.0454
$ws.Range("J136").Value = 90915860
$ws.Range("K136").Value = 14691.1362
$ws.Range("L136").Value = 272747580
$ws.Range("M136").Value = -12141.1362
$ws.Range("N136").Value = -272752680
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 22871.3
$ws.Range("J96").Value = 22871.3
$ws.Range("L96").Value = 22871.3
$ws.Range("N96").Value = -25617.3
$ws.Range("H122").Value = 3033
$ws.Range("I122").Value = 3040.8333
$ws.Range("K122").Value = 9122.499899999999
$ws.Range("M122").Value = -6672.499899999999
$ws.Range("H126").Value = 2213.7222
$ws.Range("I126").Value = 1996.2858
$ws.Range("J126").Value = 2974.75
$ws.Range("K126").Value = 5988.857400000001
$ws.Range("L126").Value = 8924.25
$ws.Range("M126").Value = -3518.857400000001
$ws.Range("N126").Value = -13864.25
$ws.Range("H132").Value = 2083.0881
$ws.Range("I132").Value = 1540.8148
$ws.Range("J132").Value = 4174.7144
$ws.Range("K132").Value = 4622.4444
$ws.Range("L132").Value = 12524.1432
$ws.Range("M132").Value = -2092.4444
$ws.Range("N132").Value = -17584.1432
$ws.Range("H136").Value = 6980.643
$ws.Range("I136").Value = 6133
$ws.Range("J136").Value = 18000
$ws.Range("K136").Value = 18399
$ws.Range("L136").Value = 54000
$ws.Range("M136").Value = -15849
$ws.Range("N136").Value = -59100
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
